# Actualizacion 11 de febrero de 2024 - Lap HP
# Se actualiza el repositorio con varios materiales.
#
# Adds week 6/7/8 attendance columns + a "Total" column to "Faltas",
# and adds P5/P6/P7 grading columns + "Puntaje"/"Calificacion" totals
# to "Concentrado".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Faltas")
$ws2 = $wb.Worksheets.Item("Concentrado")

# ---------------------------------------------------------------
# Sheet "Faltas": three new weekly date columns (J:L) + Total (M)
# ---------------------------------------------------------------

$ws1.Range("J1").Value = 45329
$ws1.Range("K1").Value = 45336
$ws1.Range("L1").Value = 45343

# Copy the date format (centered "d-mmm") from E1 onto the new date cells.
$ws1.Range("E1").Copy() | Out-Null
$ws1.Range("J1:L1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws1.Range("M1").Value = "Total"

# New attendance marks.
$ws1.Range("I3").Value = 1
$ws1.Range("I5").Value = 1
$ws1.Range("H7").Value = 1
$ws1.Range("I7").Value = 1
$ws1.Range("F8").Value = 1
$ws1.Range("H8").Value = 1
$ws1.Range("H11").Value = 1
$ws1.Range("H12").Value = 1
$ws1.Range("G13").Value = 1
$ws1.Range("H13").Value = 1
$ws1.Range("I13").Value = 1

# Totals column.
$ws1.Range("M2:M13").FormulaR1C1 = "=SUM(RC[-8]:RC[-1])"

# Apply the centered plain style (same as column E) to the new Total
# column and header cell.
$ws1.Range("E2").Copy() | Out-Null
$ws1.Range("M2:M13").PasteSpecial(-4123) | Out-Null
$excel.CutCopyMode = 0

$ws1.Range("M1").HorizontalAlignment = -4108

# Column widths (approximate best-fit widths for the new layout).
$ws1.Range("E1:I13").Columns.AutoFit() | Out-Null
$ws1.Columns("E:I").ColumnWidth = 6.86
$ws1.Columns("J:L").ColumnWidth = 6.43
$ws1.Columns("M:M").ColumnWidth = 11.43

# ---------------------------------------------------------------
# Sheet "Concentrado": P5 / P6 / P7 grading columns + Puntaje/Calificacion
# ---------------------------------------------------------------

$ws2.Range("G1").Value = "P5_Marco_Teórico"
$ws2.Range("H1").Value = "P5_Montaje"
$ws2.Range("I1").Value = "P5_Reporte"
$ws2.Range("J1").Value = "P6_Encuadre"
$ws2.Range("K1").Value = "P6_Marco_Teórico"
$ws2.Range("L1").Value = "P6_Montaje"
$ws2.Range("M1").Value = "P6_Reporte"
$ws2.Range("N1").Value = "P7_Encuadre"
$ws2.Range("O1").Value = "P7_Montaje"
$ws2.Range("P1").Value = "Puntaje"
$ws2.Range("Q1").Value = "Calificación"

# Match the header style used by the existing E1/F1 header cells.
$ws2.Range("F1").Copy() | Out-Null
$ws2.Range("G1:Q1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 2
$ws2.Range("E2").Value = 0
$ws2.Range("G2").Value = 0
$ws2.Range("H2").Value = 5
$ws2.Range("I2").Value = 0
$ws2.Range("J2").Value = 5
$ws2.Range("K2").Value = 0
$ws2.Range("L2").Value = 5
$ws2.Range("M2").Value = 0
$ws2.Range("N2").Value = 5

# Row 3
$ws2.Range("E3").Value = 0
$ws2.Range("G3").Value = 0
$ws2.Range("H3").Value = 5
$ws2.Range("I3").Value = 0
$ws2.Range("J3").Value = 5
$ws2.Range("K3").Value = 2.5
$ws2.Range("L3").Value = 5
$ws2.Range("M3").Value = 0
$ws2.Range("N3").Value = 0

# Row 4
$ws2.Range("E4").Value = 4
$ws2.Range("G4").Value = 0
$ws2.Range("H4").Value = 5
$ws2.Range("I4").Value = 0
$ws2.Range("J4").Value = 5
$ws2.Range("K4").Value = 0
$ws2.Range("L4").Value = 5
$ws2.Range("M4").Value = 0
$ws2.Range("N4").Value = 5

# Row 5
$ws2.Range("E5").Value = 4.5
$ws2.Range("G5").Value = 0
$ws2.Range("H5").Value = 5
$ws2.Range("I5").Value = 0
$ws2.Range("J5").Value = 5
$ws2.Range("K5").Value = 5
$ws2.Range("L5").Value = 5
$ws2.Range("M5").Value = 5
$ws2.Range("N5").Value = 0

# Row 6
$ws2.Range("E6").Value = 4
$ws2.Range("G6").Value = 5
$ws2.Range("H6").Value = 5
$ws2.Range("I6").Value = 5
$ws2.Range("J6").Value = 5
$ws2.Range("K6").Value = 0
$ws2.Range("L6").Value = 5
$ws2.Range("M6").Value = 0
$ws2.Range("N6").Value = 5

# Row 7
$ws2.Range("E7").Value = 4
$ws2.Range("G7").Value = 5
$ws2.Range("H7").Value = 5
$ws2.Range("I7").Value = 3
$ws2.Range("J7").Value = 5
$ws2.Range("K7").Value = 4
$ws2.Range("L7").Value = 0
$ws2.Range("M7").Value = 0
$ws2.Range("N7").Value = 0

# Row 8
$ws2.Range("E8").Value = 0
$ws2.Range("G8").Value = 0
$ws2.Range("H8").Value = 0
$ws2.Range("I8").Value = 0
$ws2.Range("J8").Value = 5
$ws2.Range("K8").Value = 0
$ws2.Range("L8").Value = 0
$ws2.Range("M8").Value = 0
$ws2.Range("N8").Value = 5

# Row 9
$ws2.Range("E9").Value = 5
$ws2.Range("G9").Value = 5
$ws2.Range("H9").Value = 5
$ws2.Range("I9").Value = 5
$ws2.Range("J9").Value = 5
$ws2.Range("K9").Value = 5
$ws2.Range("L9").Value = 5
$ws2.Range("M9").Value = 5
$ws2.Range("N9").Value = 5

# Row 10
$ws2.Range("E10").Value = 0
$ws2.Range("G10").Value = 0
$ws2.Range("H10").Value = 5
$ws2.Range("I10").Value = 0
$ws2.Range("J10").Value = 5
$ws2.Range("K10").Value = 0
$ws2.Range("L10").Value = 5
$ws2.Range("M10").Value = 0
$ws2.Range("N10").Value = 5

# Row 11
$ws2.Range("E11").Value = 0
$ws2.Range("G11").Value = 0
$ws2.Range("H11").Value = 5
$ws2.Range("I11").Value = 0
$ws2.Range("J11").Value = 5
$ws2.Range("K11").Value = 0
$ws2.Range("L11").Value = 0
$ws2.Range("M11").Value = 0
$ws2.Range("N11").Value = 5

# Row 12
$ws2.Range("E12").Value = 2.5
$ws2.Range("G12").Value = 5
$ws2.Range("H12").Value = 5
$ws2.Range("I12").Value = 0
$ws2.Range("J12").Value = 5
$ws2.Range("K12").Value = 0
$ws2.Range("L12").Value = 0
$ws2.Range("M12").Value = 0
$ws2.Range("N12").Value = 5

# Row 13
$ws2.Range("E13").Value = 0
$ws2.Range("G13").Value = 0
$ws2.Range("H13").Value = 5
$ws2.Range("I13").Value = 0
$ws2.Range("J13").Value = 0
$ws2.Range("K13").Value = 0
$ws2.Range("L13").Value = 0
$ws2.Range("M13").Value = 0
$ws2.Range("N13").Value = 0

# Apply the centered numeric style (same as column F) to every new data
# cell in E:N.
$ws2.Range("F2").Copy() | Out-Null
$ws2.Range("E2:N13").PasteSpecial(-4123) | Out-Null
$excel.CutCopyMode = 0

# Puntaje column (sum of the weighted grading columns).
$ws2.Range("P2:P13").FormulaR1C1 = "=SUM(RC[-11]:RC[-2])"
$ws2.Range("F2").Copy() | Out-Null
$ws2.Range("P2:P13").PasteSpecial(-4123) | Out-Null
$excel.CutCopyMode = 0

# Column widths for the new columns (approximate best-fit widths).
$ws2.Columns("G:G").ColumnWidth = 17.14
$ws2.Columns("H:H").ColumnWidth = 11.14
$ws2.Columns("I:I").ColumnWidth = 10.71
$ws2.Columns("J:J").ColumnWidth = 12.43
$ws2.Columns("K:K").ColumnWidth = 17.14
$ws2.Columns("L:L").ColumnWidth = 11.14
$ws2.Columns("M:M").ColumnWidth = 10.71
$ws2.Columns("N:N").ColumnWidth = 12.43
$ws2.Columns("O:O").ColumnWidth = 11.14
$ws2.Columns("P:P").ColumnWidth = 7.71
$ws2.Columns("Q:Q").ColumnWidth = 11.86

# ---------------------------------------------------------------
# Window / selection state matching the saved workbook.
# ---------------------------------------------------------------

$ws2.Activate()
$ws2.Range("M8").Select() | Out-Null

$ws1.Activate()
$ws1.Range("O11").Select() | Out-Null
